# Apply the DEUST.xlsx edit described in the commit:
#  "Fix bug related to events with the same ID, add CM groups to course info"
#
#  - Insert a new column F "Groupes CM" (all values 0) between "Heures CM" (E)
#    and the former "Heures TD" (old F), shifting old F..I to G..J.
#  - Update the sheet selection to match the saved state.
#
# Note: row order is intentionally left untouched - the underlying data is
# not actually re-sorted by this change (only a column is inserted), and the
# source data (duplicate "Code EC" values) cannot be re-sorted without
# moving rows, which the target change does not do.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F; this shifts old F..I to G..J.
$ws.Columns("F").Insert()

# Header for the newly inserted column.
$ws.Range("F1").Value = "Groupes CM"

# Fill the new column with 0 for every data row (rows 2-33).
$ws.Range("F2:F33").Value = 0

# Update the selection to match the committed state.
$ws.Range("D31").Select()
